$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-09 Monday", 2) | Out-Null
$d.Content.Find.Execute("26+36=", $true, $false, $false, $false, $false, $true, 1, $false, "85-16=", 2) | Out-Null
$d.Content.Find.Execute("64+4=", $true, $false, $false, $false, $false, $true, 1, $false, "61+19=", 2) | Out-Null
$d.Content.Find.Execute("44+54=", $true, $false, $false, $false, $false, $true, 1, $false, "41-20=", 2) | Out-Null
$d.Content.Find.Execute("98-12=", $true, $false, $false, $false, $false, $true, 1, $false, "93-62=", 2) | Out-Null
$d.Content.Find.Execute("91-27=", $true, $false, $false, $false, $false, $true, 1, $false, "53-17=", 2) | Out-Null
$d.Content.Find.Execute("44+24=", $true, $false, $false, $false, $false, $true, 1, $false, "81-55=", 2) | Out-Null
$d.Content.Find.Execute("44-32=", $true, $false, $false, $false, $false, $true, 1, $false, "43+8=", 2) | Out-Null
$d.Content.Find.Execute("37+18=", $true, $false, $false, $false, $false, $true, 1, $false, "91-59=", 2) | Out-Null
$d.Content.Find.Execute("66+25=", $true, $false, $false, $false, $false, $true, 1, $false, "17+79=", 2) | Out-Null
$d.Content.Find.Execute("4+28=", $true, $false, $false, $false, $false, $true, 1, $false, "59-23=", 2) | Out-Null
$d.Content.Find.Execute("40+51=", $true, $false, $false, $false, $false, $true, 1, $false, "73-48=", 2) | Out-Null
$d.Content.Find.Execute("13+28=", $true, $false, $false, $false, $false, $true, 1, $false, "59-43=", 2) | Out-Null
$d.Content.Find.Execute("66+3=", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=", 2) | Out-Null
$d.Content.Find.Execute("21+29=", $true, $false, $false, $false, $false, $true, 1, $false, "11+30=", 2) | Out-Null
$d.Content.Find.Execute("66-65=", $true, $false, $false, $false, $false, $true, 1, $false, "65+18=", 2) | Out-Null
$d.Content.Find.Execute("42-42=", $true, $false, $false, $false, $false, $true, 1, $false, "7+32=", 2) | Out-Null
$d.Content.Find.Execute("92-6=", $true, $false, $false, $false, $false, $true, 1, $false, "67-19=", 2) | Out-Null
$d.Content.Find.Execute("1+77=", $true, $false, $false, $false, $false, $true, 1, $false, "86-74=", 2) | Out-Null
$d.Content.Find.Execute("38+2=", $true, $false, $false, $false, $false, $true, 1, $false, "51+42=", 2) | Out-Null
$d.Content.Find.Execute("8+75=", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=", 2) | Out-Null
$d.Content.Find.Execute("28-10=", $true, $false, $false, $false, $false, $true, 1, $false, "81-72=", 2) | Out-Null
$d.Content.Find.Execute("98-9=", $true, $false, $false, $false, $false, $true, 1, $false, "32-16=", 2) | Out-Null
$d.Content.Find.Execute("36+7=", $true, $false, $false, $false, $false, $true, 1, $false, "45+12=", 2) | Out-Null
$d.Content.Find.Execute("33-22=", $true, $false, $false, $false, $false, $true, 1, $false, "90-29=", 2) | Out-Null
$d.Content.Find.Execute("53+44=", $true, $false, $false, $false, $false, $true, 1, $false, "44+48=", 2) | Out-Null
$d.Content.Find.Execute("2+61=", $true, $false, $false, $false, $false, $true, 1, $false, "93-43=", 2) | Out-Null
$d.Content.Find.Execute("97-76=", $true, $false, $false, $false, $false, $true, 1, $false, "10+29=", 2) | Out-Null
$d.Content.Find.Execute("5+93=", $true, $false, $false, $false, $false, $true, 1, $false, "79-4=", 2) | Out-Null
$d.Content.Find.Execute("71-20=", $true, $false, $false, $false, $false, $true, 1, $false, "12+44=", 2) | Out-Null
$d.Content.Find.Execute("8+47=", $true, $false, $false, $false, $false, $true, 1, $false, "43+4=", 2) | Out-Null
$d.Content.Find.Execute("89-30=", $true, $false, $false, $false, $false, $true, 1, $false, "27+4=", 2) | Out-Null
$d.Content.Find.Execute("38+26=", $true, $false, $false, $false, $false, $true, 1, $false, "9+71=", 2) | Out-Null
$d.Content.Find.Execute("54+33=", $true, $false, $false, $false, $false, $true, 1, $false, "38-34=", 2) | Out-Null
$d.Content.Find.Execute("25-15=", $true, $false, $false, $false, $false, $true, 1, $false, "0+35=", 2) | Out-Null
$d.Content.Find.Execute("20+54=", $true, $false, $false, $false, $false, $true, 1, $false, "22+0=", 2) | Out-Null
$d.Content.Find.Execute("62-6=", $true, $false, $false, $false, $false, $true, 1, $false, "90-1=", 2) | Out-Null
$d.Content.Find.Execute("37-16=", $true, $false, $false, $false, $false, $true, 1, $false, "0+57=", 2) | Out-Null
$d.Content.Find.Execute("31-3=", $true, $false, $false, $false, $false, $true, 1, $false, "9+81=", 2) | Out-Null
$d.Content.Find.Execute("20+75=", $true, $false, $false, $false, $false, $true, 1, $false, "25-8=", 2) | Out-Null
$d.Content.Find.Execute("13+17=", $true, $false, $false, $false, $false, $true, 1, $false, "68-53=", 2) | Out-Null
$d.Content.Find.Execute("18+23=", $true, $false, $false, $false, $false, $true, 1, $false, "4+0=", 2) | Out-Null
$d.Content.Find.Execute("81-46=", $true, $false, $false, $false, $false, $true, 1, $false, "33+42=", 2) | Out-Null
$d.Content.Find.Execute("76-69=", $true, $false, $false, $false, $false, $true, 1, $false, "65-61=", 2) | Out-Null
$d.Content.Find.Execute("1+32=", $true, $false, $false, $false, $false, $true, 1, $false, "41-10=", 2) | Out-Null
$d.Content.Find.Execute("51+41=", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=", 2) | Out-Null
$d.Content.Find.Execute("39-15=", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=", 2) | Out-Null
$d.Content.Find.Execute("81-61=", $true, $false, $false, $false, $false, $true, 1, $false, "73-28=", 2) | Out-Null
$d.Content.Find.Execute("10+45=", $true, $false, $false, $false, $false, $true, 1, $false, "96-38=", 2) | Out-Null
$d.Content.Find.Execute("48+41=", $true, $false, $false, $false, $false, $true, 1, $false, "6+55=", 2) | Out-Null
$d.Content.Find.Execute("58+17=", $true, $false, $false, $false, $false, $true, 1, $false, "42+39=", 2) | Out-Null
$d.Content.Find.Execute("41+0=", $true, $false, $false, $false, $false, $true, 1, $false, "56-46=", 2) | Out-Null
$d.Content.Find.Execute("84-62=", $true, $false, $false, $false, $false, $true, 1, $false, "76+17=", 2) | Out-Null
$d.Content.Find.Execute("51-0=", $true, $false, $false, $false, $false, $true, 1, $false, "5+46=", 2) | Out-Null
$d.Content.Find.Execute("40-5=", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=", 2) | Out-Null
$d.Content.Find.Execute("9+65=", $true, $false, $false, $false, $false, $true, 1, $false, "85-63=", 2) | Out-Null
$d.Content.Find.Execute("10+51=", $true, $false, $false, $false, $false, $true, 1, $false, "90-85=", 2) | Out-Null
$d.Content.Find.Execute("36+28=", $true, $false, $false, $false, $false, $true, 1, $false, "24-15=", 2) | Out-Null
$d.Content.Find.Execute("83+7=", $true, $false, $false, $false, $false, $true, 1, $false, "52-44=", 2) | Out-Null
$d.Content.Find.Execute("11+49=", $true, $false, $false, $false, $false, $true, 1, $false, "94-74=", 2) | Out-Null
$d.Content.Find.Execute("58+29=", $true, $false, $false, $false, $false, $true, 1, $false, "56-35=", 2) | Out-Null
$d.Content.Find.Execute("92-13=", $true, $false, $false, $false, $false, $true, 1, $false, "16+23=", 2) | Out-Null
$d.Content.Find.Execute("91-66=", $true, $false, $false, $false, $false, $true, 1, $false, "50-6=", 2) | Out-Null
$d.Content.Find.Execute("41-25=", $true, $false, $false, $false, $false, $true, 1, $false, "31+67=", 2) | Out-Null
$d.Content.Find.Execute("81+11=", $true, $false, $false, $false, $false, $true, 1, $false, "35+18=", 2) | Out-Null
$d.Content.Find.Execute("8+25=", $true, $false, $false, $false, $false, $true, 1, $false, "43+11=", 2) | Out-Null
$d.Content.Find.Execute("23+18=", $true, $false, $false, $false, $false, $true, 1, $false, "91-68=", 2) | Out-Null
$d.Content.Find.Execute("82-48=", $true, $false, $false, $false, $false, $true, 1, $false, "81+14=", 2) | Out-Null
$d.Content.Find.Execute("86-52=", $true, $false, $false, $false, $false, $true, 1, $false, "95+1=", 2) | Out-Null
$d.Content.Find.Execute("9+39=", $true, $false, $false, $false, $false, $true, 1, $false, "87+11=", 2) | Out-Null
$d.Content.Find.Execute("94-15=", $true, $false, $false, $false, $false, $true, 1, $false, "75-32=", 2) | Out-Null
$d.Content.Find.Execute("90+0=", $true, $false, $false, $false, $false, $true, 1, $false, "89-19=", 2) | Out-Null
$d.Content.Find.Execute("95-25=", $true, $false, $false, $false, $false, $true, 1, $false, "74+21=", 2) | Out-Null
$d.Content.Find.Execute("63-44=", $true, $false, $false, $false, $false, $true, 1, $false, "19+8=", 2) | Out-Null
$d.Content.Find.Execute("72+12=", $true, $false, $false, $false, $false, $true, 1, $false, "65-35=", 2) | Out-Null
$d.Content.Find.Execute("46+25=", $true, $false, $false, $false, $false, $true, 1, $false, "17+2=", 2) | Out-Null
$d.Content.Find.Execute("88-48=", $true, $false, $false, $false, $false, $true, 1, $false, "78-15=", 2) | Out-Null
$d.Content.Find.Execute("54-44=", $true, $false, $false, $false, $false, $true, 1, $false, "63-45=", 2) | Out-Null
$d.Content.Find.Execute("85-41=", $true, $false, $false, $false, $false, $true, 1, $false, "60-48=", 2) | Out-Null
$d.Content.Find.Execute("83+2=", $true, $false, $false, $false, $false, $true, 1, $false, "15+11=", 2) | Out-Null
$d.Content.Find.Execute("22-11=", $true, $false, $false, $false, $false, $true, 1, $false, "89+8=", 2) | Out-Null
$d.Content.Find.Execute("51-41=", $true, $false, $false, $false, $false, $true, 1, $false, "52-50=", 2) | Out-Null
$d.Content.Find.Execute("42+45=", $true, $false, $false, $false, $false, $true, 1, $false, "26+29=", 2) | Out-Null
$d.Content.Find.Execute("18+64=", $true, $false, $false, $false, $false, $true, 1, $false, "20+65=", 2) | Out-Null
$d.Content.Find.Execute("93-50=", $true, $false, $false, $false, $false, $true, 1, $false, "44+52=", 2) | Out-Null
$d.Content.Find.Execute("55+40=", $true, $false, $false, $false, $false, $true, 1, $false, "11-1=", 2) | Out-Null
$d.Content.Find.Execute("57-52=", $true, $false, $false, $false, $false, $true, 1, $false, "24+12=", 2) | Out-Null
$d.Content.Find.Execute("66-33=", $true, $false, $false, $false, $false, $true, 1, $false, "72+13=", 2) | Out-Null
$d.Content.Find.Execute("67-13=", $true, $false, $false, $false, $false, $true, 1, $false, "40-24=", 2) | Out-Null
$d.Content.Find.Execute("60+38=", $true, $false, $false, $false, $false, $true, 1, $false, "7+75=", 2) | Out-Null
$d.Content.Find.Execute("17+18=", $true, $false, $false, $false, $false, $true, 1, $false, "66+10=", 2) | Out-Null
$d.Content.Find.Execute("99-97=", $true, $false, $false, $false, $false, $true, 1, $false, "78-50=", 2) | Out-Null
$d.Content.Find.Execute("62-51=", $true, $false, $false, $false, $false, $true, 1, $false, "6+93=", 2) | Out-Null
$d.Content.Find.Execute("4+49=", $true, $false, $false, $false, $false, $true, 1, $false, "21+61=", 2) | Out-Null
$d.Content.Find.Execute("26+8=", $true, $false, $false, $false, $false, $true, 1, $false, "39+31=", 2) | Out-Null
$d.Content.Find.Execute("0+62=", $true, $false, $false, $false, $false, $true, 1, $false, "44-41=", 2) | Out-Null
$d.Content.Find.Execute("80-72=", $true, $false, $false, $false, $false, $true, 1, $false, "24+36=", 2) | Out-Null
$d.Content.Find.Execute("15-10=", $true, $false, $false, $false, $false, $true, 1, $false, "80-41=", 2) | Out-Null
$d.Content.Find.Execute("67+13=", $true, $false, $false, $false, $false, $true, 1, $false, "59-45=", 2) | Out-Null
$d.Content.Find.Execute("80-58=", $true, $false, $false, $false, $false, $true, 1, $false, "23+12=", 2) | Out-Null
$d.Content.Find.Execute("87+0=", $true, $false, $false, $false, $false, $true, 1, $false, "29+26=", 2) | Out-Null
